$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# --- Row 11 ---
$ws.Range("F11").Value = "Hecho"
$ws.Range("N11").Value = 1.5
$ws.Range("AU11").Value = $null
$ws.Range("AX11").Value = $null

# --- Row 21 ---
$ws.Range("F21").Value = "Hecho"
$ws.Range("W21").Value = 0.5

# --- Row 24 ---
$ws.Range("Q24").Value = 2
$ws.Range("W24").Value = 0.5
$ws.Range("AU24").Value = $null

# --- Row 30 ---
$ws.Range("F30").Value = "Hecho"
$ws.Range("AF30").Value = $null
$ws.Range("AI30").Value = 1
$ws.Range("AU30").Value = $null

# --- Row 33 ---
$ws.Range("AC33").Value = 2
$ws.Range("AU33").Value = $null
$ws.Range("AX33").Value = $null

# --- Row 36 ---
$ws.Range("F36").Value = "En proceso"

# --- Row 38 ---
$ws.Range("F38").Value = "Hecho"
$ws.Range("W38").Value = 0.5
$ws.Range("AL38").Value = $null

# --- Row 39 ---
$ws.Range("F39").Value = "Hecho"
$ws.Range("AI39").Value = 0.5
$ws.Range("AR39").Value = $null

# --- Row 40 ---
$ws.Range("F40").Value = "Hecho"
$ws.Range("W40").Value = 0.75
$ws.Range("AC40").Value = $null
$ws.Range("AU40").Value = $null

# --- Row 41 ---
$ws.Range("F41").Value = "Hecho"
$ws.Range("AC41").Value = 1.5
